# Replace the semicolon separator used in credible-interval strings like
# "(-1,550; 4,255)" with the word "to", i.e. "(-1,550 to 4,255)".
# This applies uniformly to every text cell in the used range that
# contains the "; " substring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            if ($val.Contains("; ")) {
                $cell.Value2 = $val.Replace("; ", " to ")
            }
        }
    }
}
